# Fruta / hortaliza, semanal
# Two new weekly price rows were added at the top of the "Durazno"
# dataset (Vega Modelo de Temuco). Inserting at row 324 shifts every
# existing record down by two rows (324->326 ... 418->420), which is
# exactly what the canonical diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 324, pushing the rest of
# the table (old rows 324..418) down to 326..420.
$ws.Rows.Item(324).Insert()
$ws.Rows.Item(324).Insert()

# --- New row 324: Elegant Lady ---
$ws.Cells.Item(324, 1).Value = 10
$ws.Cells.Item(324, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(324, 3).Value = "La Araucanía"
$ws.Cells.Item(324, 4).Value = 44988
$ws.Cells.Item(324, 5).Value = 9
$ws.Cells.Item(324, 6).Value = "Fruta"
$ws.Cells.Item(324, 7).Value = 100103
$ws.Cells.Item(324, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(324, 9).Value = 100103004
$ws.Cells.Item(324, 10).Value = "Durazno"
$ws.Cells.Item(324, 11).Value = "Elegant Lady"
$ws.Cells.Item(324, 12).Value = "Primera"
$ws.Cells.Item(324, 13).Value = 125
$ws.Cells.Item(324, 14).Value = 21000
$ws.Cells.Item(324, 15).Value = 21000
$ws.Cells.Item(324, 16).Value = 21000
$ws.Cells.Item(324, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(324, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(324, 19).Value = 1167
$ws.Cells.Item(324, 20).Value = 18

# --- New row 325: September Sun ---
$ws.Cells.Item(325, 1).Value = 10
$ws.Cells.Item(325, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(325, 3).Value = "La Araucanía"
$ws.Cells.Item(325, 4).Value = 44988
$ws.Cells.Item(325, 5).Value = 9
$ws.Cells.Item(325, 6).Value = "Fruta"
$ws.Cells.Item(325, 7).Value = 100103
$ws.Cells.Item(325, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(325, 9).Value = 100103004
$ws.Cells.Item(325, 10).Value = "Durazno"
$ws.Cells.Item(325, 11).Value = "September Sun"
$ws.Cells.Item(325, 12).Value = "Primera"
$ws.Cells.Item(325, 13).Value = 110
$ws.Cells.Item(325, 14).Value = 21000
$ws.Cells.Item(325, 15).Value = 21000
$ws.Cells.Item(325, 16).Value = 21000
$ws.Cells.Item(325, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(325, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(325, 19).Value = 1167
$ws.Cells.Item(325, 20).Value = 18
